$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "https://github.com/evansagge/mongoid-rspec"
$ws.Range("A2").Value = "ROR"
$ws.Range("B2").Value = "shulda matcher not supporting for writing rspec model. this model using mongodb databases. Use the gem called mongoid-rspec. "

$ws.Range("B2").Select()
